$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "Login" case, drop the "senha" half from B2/C2 ---
$ws.Range("B2").Value = "login = bobbmarcio"

$ws.Range("C2").Value = "login = null"
$ws.Range("C2").Characters(9, 4).Font.Italic = $true

$ws.Range("D2").Value = "CE_Login"

# --- Row 3: was "Criar Login" case, now becomes the "Senha" case ---
$ws.Range("A3").Value = "Senha"
$ws.Range("B3").Value = "senha = 12345678"

$ws.Range("C3").Value = "senha = null"
$ws.Range("C3").Characters(9, 4).Font.Italic = $true

$ws.Range("D3").Value = "CE_Senha"

# --- Row 4: new "Nome do campeonato" case ---
$ws.Range("A4").Value = "Nome do campeonato"
$ws.Range("B4").Value = "nome = Campeonato Brasileiro"

$ws.Range("C4").Value = "nome = null"
$ws.Range("C4").Characters(8, 4).Font.Italic = $true

$ws.Range("D4").Value = "CE_NomeCampeonato"

# --- Row 5: new "Quantidade de times" case ---
$ws.Range("A5").Value = "Quantidade de times"
$ws.Range("B5").Value = "quantidade = 3..."
$ws.Range("C5").Value = "quantidade = ..., -2, -1, 0, 1, 2"
$ws.Range("D5").Value = "CE_QuantidadeTimes"

# --- Row 6: new "Nome do time" case (reuses the same rich "nome = null" text as row 4) ---
$ws.Range("A6").Value = "Nome do time"
$ws.Range("B6").Value = "nome = Corinthinas"

$ws.Range("C4").Copy()
$ws.Range("C6").PasteSpecial()
$excel.CutCopyMode = $false

$ws.Range("D6").Value = "CE_NomeTime"

# --- Header row shading ---
$ws.Range("A1:D1").Interior.Color = 12566463

# --- Column widths (auto fit to new, wider content in columns A/B) ---
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()

# --- Final selection ---
$ws.Range("A7").Select()
